$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; existing rows 16-45 shift down to 17-46.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the new weekly data point
# (same Mercado/Región/Categoría/etc. metadata as the surrounding rows,
# new Fecha + price figures).
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = "Femacal de La Calera"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44880
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 300000000
$ws.Range("G16").Value = "Espárragos"
$ws.Range("H16").Value = "Verde"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1300
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = 1355
$ws.Range("N16").Value = "$/kilo"
$ws.Range("O16").Value = "Provincia de Quillota"
$ws.Range("P16").Value = 1355
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"
